# Adding the changes we made on may 9th
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new blank rows starting at row 3 (kept away from the header row
# so the blank rows don't pick up the header's formatting). This pushes
# the old row 2 down to row 11, and old rows 3-21 down to rows 12-30.
$ws.Range("A3:C11").EntireRow.Insert()

# The old row 2 (still sitting at row 2) now needs to move down to row 11
# to join the rest of the shifted data.
$ws.Range("A2:C2").Cut($ws.Range("A11:C11"))

# Fill rows 2-10 with the new data.
$newTop = @(
    @(0.0500909499824047, 0.0342084541916847, 0.0232128798961639),
    @(0.0068722339347004, 0.0074830991216003, 0.0390953756868839),
    @(-0.0138971842825412, 0.0290160998702049, 0.1440114825963974),
    @(0.0100792767480015, -0.022754730656743, 0.0288633834570646),
    @(-0.0282525178045034, -0.020616702735424, -0.0172569435089826),
    @(-0.0167987942695617, -0.0216857157647609, 0),
    @(-0.0178678091615438, 0.0054977871477603, 0.0299323964864015),
    @(-0.0229074470698833, 0.00534507073462, 0.0030543261673301),
    @(-0.0394008085131645, 0.0178678091615438, 0.011148290708661)
)

$r = 2
foreach ($row in $newTop) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Append one new row of data at the end (row 31).
$ws.Cells.Item(31, 1).Value = -0.3419318199157715
$ws.Cells.Item(31, 2).Value = 2.081828832626343
$ws.Cells.Item(31, 3).Value = -0.5186246037483215
